$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the duplicated state+country text fragments in column P (address line),
# normalizing them all to the "... usa" suffix form.
$ws.Range("P5").Value = "123 Connecticut CT usa"
$ws.Range("P6").Value = "123 Delaware DE usa"
$ws.Range("P8").Value = "123 hawai hi usa"
$ws.Range("P9").Value = "123 Iowa IA usa"
$ws.Range("P10").Value = "123 Idaho ID usa"
$ws.Range("P11").Value = "123 Illinois IL usa"
$ws.Range("P12").Value = "123 Indiana IN usa"
$ws.Range("P13").Value = "123 Kansas KS usa"
$ws.Range("P14").Value = "123 Massachusetts MA usa"
$ws.Range("P15").Value = "123 Maryland MD usa"
$ws.Range("P16").Value = "123 Maine ME usa"
$ws.Range("P18").Value = "123 Minnesota MN usa"
$ws.Range("P20").Value = "123 Mississippi MS usa"
$ws.Range("P21").Value = "123 Montana MT usa"
$ws.Range("P22").Value = "123 North Dakota ND usa"
$ws.Range("P23").Value = "123 Nebraska NE usa"
$ws.Range("P24").Value = "123 New Hampshire NH usa"
$ws.Range("P25").Value = "123 New Mexico NM usa"
$ws.Range("P26").Value = "123 Nevada NV usa"
$ws.Range("P27").Value = "123 New York NY usa"
$ws.Range("P28").Value = "123 Ohio OH usa"
$ws.Range("P29").Value = "123 Oklahoma OK usa"
$ws.Range("P31").Value = "123 Pennsylvania PA usa"
$ws.Range("P33").Value = "123 South Carolina SC usa"
$ws.Range("P34").Value = "123 South Dakota SD usa"
$ws.Range("P35").Value = "123 Tennessee TN usa"
$ws.Range("P36").Value = "123 Texas TX usa"
$ws.Range("P37").Value = "123 Utah UT usa"
$ws.Range("P38").Value = "123 Virginia VA usa"
$ws.Range("P39").Value = "123 Vermont VT usa"
$ws.Range("P40").Value = "123 Washington WA usa"
$ws.Range("P41").Value = "123 Wisconsin WI usa"

# Leave the cursor where the editor last clicked after making these edits.
$ws.Range("N10").Select()
